$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking text such as "236.31" or "92.242.79".
# Force the cells to Text format first so Excel keeps these values as strings
# instead of silently converting them to numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '92.242.79'
$ws.Range("E2").Value = '  +0.83%  '

$ws.Range("D3").Value = '3.092.54'
$ws.Range("E3").Value = '  -2.03%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '236.31'
$ws.Range("E5").Value = '  -1.12%  '

$ws.Range("D6").Value = '609.85'
$ws.Range("E6").Value = '  -1.83%  '

$ws.Range("E7").Value = '  -3.66%  '

$ws.Range("D8").Value = '0.390'
$ws.Range("E8").Value = '  +4.04%  '

$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.10%  '

$ws.Range("D10").Value = '3.086.94'
$ws.Range("E10").Value = '  -2.16%  '

$ws.Range("D11").Value = '0.732'
$ws.Range("E11").Value = '  -1.66%  '

$ws.Range("D12").Value = '0.200'
$ws.Range("E12").Value = '  -1.54%  '

$ws.Range("D13").Value = '0.0000247'
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("D14").Value = '92.350.78'
$ws.Range("E14").Value = '  +1.33%  '

$ws.Range("D15").Value = '34.07'
$ws.Range("E15").Value = '  -3.96%  '

$ws.Range("D16").Value = '5.42'
$ws.Range("E16").Value = '  -3.00%  '

$ws.Range("D17").Value = '3.670.53'
$ws.Range("E17").Value = '  -2.14%  '

$ws.Range("D18").Value = '3.100.95'
$ws.Range("E18").Value = '  -1.95%  '

$ws.Range("D19").Value = '3.76'
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").Value = '14.60'
$ws.Range("E20").Value = '  -4.60%  '

$ws.Range("D21").Value = '5.73'
$ws.Range("E21").Value = '  -3.96%  '

$ws.Range("D22").Value = '9.31'
$ws.Range("E22").Value = '  +1.42%  '

$ws.Range("D23").Value = '442.47'
$ws.Range("E23").Value = '  -3.26%  '

$ws.Range("D24").Value = '0.0000194'
$ws.Range("E24").Value = '  -4.86%  '

$ws.Range("D25").Value = '5.69'
$ws.Range("E25").Value = '  -5.58%  '

$ws.Range("D26").Value = '85.87'
$ws.Range("E26").Value = '  -3.71%  '

$ws.Range("D27").Value = '11.61'
$ws.Range("E27").Value = '  -3.71%  '

$ws.Range("D28").Value = '3.253.90'
$ws.Range("E28").Value = '  -1.95%  '

$ws.Range("E29").Value = '  -0.22%  '

$ws.Range("D30").Value = '0.129'
$ws.Range("E30").Value = '  +0.92%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.227'
$ws.Range("E31").Value = '  -1.41%  '

$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").Value = '0.168'
$ws.Range("E32").Value = '  -2.41%  '

$ws.Range("D33").Value = '9.10'
$ws.Range("E33").Value = '  -3.05%  '

$ws.Range("E34").Value = '  +6.60%  '

$ws.Range("D35").Value = '7.86'
$ws.Range("E35").Value = '  +2.54%  '

$ws.Range("E36").Value = '  -8.16%  '

$ws.Range("D37").Value = '25.82'
$ws.Range("E37").Value = '  -2.69%  '

$ws.Range("E38").Value = '  -3.66%  '

$ws.Range("D39").Value = '3.86'
$ws.Range("E39").Value = '  +0.78%  '

$ws.Range("D40").Value = '483.97'
$ws.Range("E40").Value = '  -5.66%  '

$ws.Range("D41").Value = '23.87'
$ws.Range("E41").Value = '  +7.55%  '

$ws.Range("E42").Value = '  -5.27%  '

$ws.Range("D43").Value = '0.430'
$ws.Range("E43").Value = '  -5.14%  '

$ws.Range("D44").Value = '3.30'
$ws.Range("E44").Value = '  -4.65%  '

$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").Value = '164.82'
$ws.Range("E46").Value = '  +4.41%  '

$ws.Range("E47").Value = '  -3.89%  '

$ws.Range("D48").Value = '0.682'
$ws.Range("E48").Value = '  -4.23%  '

$ws.Range("D49").Value = '1.38'
$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").Value = '0.0332'
$ws.Range("E50").Value = '  +3.96%  '

$ws.Range("D51").Value = '43.94'
$ws.Range("E51").Value = '  -0.32%  '
